$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.244.35'
$ws.Range("D3").Value = '3.152.38'
$ws.Range("E3").Value = '  +5.88%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.43'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.78%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.144.50'
$ws.Range("E8").Value = '  +5.81%  '
$ws.Range("E9").Value = '  +3.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +21.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.76'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.52%  '
$ws.Range("E12").Value = '  +5.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000257'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +12.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.62%  '
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '3.673.54'
$ws.Range("E16").Value = '  +5.91%  '
$ws.Range("D17").Value = '64.152.24'
$ws.Range("E17").Value = '  +8.44%  '
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = '3.151.92'
$ws.Range("E19").Value = '  +5.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.733'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.37%  '
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +13.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.92%  '
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.25'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.27%  '
$ws.Range("E33").Value = '  +6.71%  '
$ws.Range("D34").Value = '0.0₃0880'
$ws.Range("E34").Value = '  +15.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +19.07%  '
$ws.Range("E36").Value = '  +7.20%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +21.41%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.15'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '51.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '449.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.26%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").Value = '2.941.93'
$ws.Range("E42").Value = '  +8.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0373'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.47%  '
$ws.Range("E44").Value = '  +14.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.114'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.27%  '
$ws.Range("E46").Value = '  +12.62%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.02'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.24'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").Value = '  +2.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.34%  '
